$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2435.85
$ws.Cells.Item(17, 10).Value = 2435.85
$ws.Cells.Item(17, 12).Value = 7307.549999999999
$ws.Cells.Item(17, 14).Value = -7643.549999999999

$ws.Cells.Item(28, 8).Value = 1178.9166
$ws.Cells.Item(28, 9).Value = 1224.7
$ws.Cells.Item(28, 10).Value = 950
$ws.Cells.Item(28, 11).Value = 1224.7
$ws.Cells.Item(28, 12).Value = 950
$ws.Cells.Item(28, 13).Value = -739.7
$ws.Cells.Item(28, 14).Value = -1920

$ws.Cells.Item(31, 8).Value = 23.5
$ws.Cells.Item(31, 9).Value = 23.5
$ws.Cells.Item(31, 11).Value = 70.5
$ws.Cells.Item(31, 13).Value = 159.5

$ws.Cells.Item(32, 8).Value = 3409.7
$ws.Cells.Item(32, 9).Value = 1888.5
$ws.Cells.Item(32, 11).Value = 1888.5
$ws.Cells.Item(32, 13).Value = -1562.5

$ws.Cells.Item(40, 8).Value = 3579.7
$ws.Cells.Item(40, 9).Value = 3250.5
$ws.Cells.Item(40, 10).Value = 3799.1667
$ws.Cells.Item(40, 11).Value = 3250.5
$ws.Cells.Item(40, 12).Value = 3799.1667
$ws.Cells.Item(40, 13).Value = -3075.5
$ws.Cells.Item(40, 14).Value = -4149.1667

$ws.Cells.Item(51, 8).Value = 10191.846
$ws.Cells.Item(51, 9).Value = 9785
$ws.Cells.Item(51, 10).Value = 10666.5
$ws.Cells.Item(51, 11).Value = 9785
$ws.Cells.Item(51, 12).Value = 10666.5
$ws.Cells.Item(51, 13).Value = -9301
$ws.Cells.Item(51, 14).Value = -11634.5

$ws.Cells.Item(55, 8).Value = 1007.1053
$ws.Cells.Item(55, 9).Value = 1629.1111
$ws.Cells.Item(55, 10).Value = 447.3
$ws.Cells.Item(55, 11).Value = 1629.1111
$ws.Cells.Item(55, 12).Value = 447.3
$ws.Cells.Item(55, 13).Value = -1415.1111
$ws.Cells.Item(55, 14).Value = -875.3

$ws.Cells.Item(113, 8).Value = 4099.7334
$ws.Cells.Item(113, 10).Value = 2499.5
$ws.Cells.Item(113, 12).Value = 2499.5
$ws.Cells.Item(113, 14).Value = -9007.5

$ws.Cells.Item(125, 8).Value = 762.5
$ws.Cells.Item(125, 10).Value = 618.125
$ws.Cells.Item(125, 12).Value = 5563.125
$ws.Cells.Item(125, 14).Value = -10483.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 3000
$ws.Cells.Item(3, 10).Value = 3000
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 14).Value = -3230

$ws.Cells.Item(22, 8).Value = 12562.25
$ws.Cells.Item(22, 9).Value = 6724.5
$ws.Cells.Item(22, 10).Value = 18400
$ws.Cells.Item(22, 11).Value = 6724.5
$ws.Cells.Item(22, 12).Value = 18400
$ws.Cells.Item(22, 13).Value = -6425.5
$ws.Cells.Item(22, 14).Value = -18998

$ws.Cells.Item(37, 8).Value = 10583.667
$ws.Cells.Item(37, 9).Value = 10583.667
$ws.Cells.Item(37, 11).Value = 10583.667
$ws.Cells.Item(37, 13).Value = -10310.667

$ws.Cells.Item(61, 8).Value = 6869
$ws.Cells.Item(61, 9).Value = 5037.5
$ws.Cells.Item(61, 11).Value = 5037.5
$ws.Cells.Item(61, 13).Value = -4825.5

$ws.Cells.Item(132, 8).Value = 6305.75
$ws.Cells.Item(132, 9).Value = 6585
$ws.Cells.Item(132, 10).Value = 5747.25
$ws.Cells.Item(132, 11).Value = 19755
$ws.Cells.Item(132, 12).Value = 17241.75
$ws.Cells.Item(132, 13).Value = -17225
$ws.Cells.Item(132, 14).Value = -22301.75

$ws.Cells.Item(136, 8).Value = 6869
$ws.Cells.Item(136, 9).Value = 5037.5
$ws.Cells.Item(136, 11).Value = 15112.5
$ws.Cells.Item(136, 13).Value = -12562.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2299.6667
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 10).Value = 2949.5
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 12).Value = 2949.5
$ws.Cells.Item(107, 13).Value = 920
$ws.Cells.Item(107, 14).Value = -6789.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 5833.3335
$ws.Cells.Item(4, 10).Value = 5833.3335
$ws.Cells.Item(4, 12).Value = 5833.3335
$ws.Cells.Item(4, 14).Value = -6057.3335

$ws.Cells.Item(7, 8).Value = 221.1
$ws.Cells.Item(7, 9).Value = 179
$ws.Cells.Item(7, 10).Value = 600
$ws.Cells.Item(7, 11).Value = 179
$ws.Cells.Item(7, 12).Value = 600
$ws.Cells.Item(7, 13).Value = -66
$ws.Cells.Item(7, 14).Value = -826

$ws.Cells.Item(10, 8).Value = 3400
$ws.Cells.Item(10, 9).Value = 231
$ws.Cells.Item(10, 10).Value = 6569
$ws.Cells.Item(10, 11).Value = 231
$ws.Cells.Item(10, 12).Value = 6569
$ws.Cells.Item(10, 13).Value = -92
$ws.Cells.Item(10, 14).Value = -6847

$ws.Cells.Item(17, 8).Value = 1650
$ws.Cells.Item(17, 9).Value = 1600
$ws.Cells.Item(17, 10).Value = 1700
$ws.Cells.Item(17, 11).Value = 1600
$ws.Cells.Item(17, 12).Value = 1700
$ws.Cells.Item(17, 13).Value = -1426
$ws.Cells.Item(17, 14).Value = -2048

$ws.Cells.Item(60, 8).Value = 29133.334
$ws.Cells.Item(60, 9).Value = 29133.334
$ws.Cells.Item(60, 11).Value = 29133.334
$ws.Cells.Item(60, 13).Value = -28622.334

$ws.Cells.Item(96, 8).Value = 15210
$ws.Cells.Item(96, 10).Value = 15210
$ws.Cells.Item(96, 12).Value = 15210
$ws.Cells.Item(96, 14).Value = -20702

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(35, 8).Value = 10000
$ws.Cells.Item(35, 9).Value = 10000
$ws.Cells.Item(35, 11).Value = 30000
$ws.Cells.Item(35, 13).Value = -29712

$ws.Cells.Item(37, 8).Value = 190000
$ws.Cells.Item(37, 10).Value = 190000
$ws.Cells.Item(37, 12).Value = 570000
$ws.Cells.Item(37, 14).Value = -570224

$ws.Cells.Item(131, 8).Value = 772
$ws.Cells.Item(131, 9).Value = 663.2857
$ws.Cells.Item(131, 10).Value = 1533
$ws.Cells.Item(131, 11).Value = 1989.8571
$ws.Cells.Item(131, 12).Value = 4599
$ws.Cells.Item(131, 13).Value = 3050.1429
$ws.Cells.Item(131, 14).Value = -14679

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2404.875
$ws.Cells.Item(40, 9).Value = 1772.3334
$ws.Cells.Item(40, 10).Value = 4302.5
$ws.Cells.Item(40, 11).Value = 1772.3334
$ws.Cells.Item(40, 12).Value = 4302.5
$ws.Cells.Item(40, 13).Value = -1636.3334
$ws.Cells.Item(40, 14).Value = -4574.5

$ws.Cells.Item(46, 8).Value = 2197.1667
$ws.Cells.Item(46, 10).Value = 3184.3333
$ws.Cells.Item(46, 12).Value = 3184.3333
$ws.Cells.Item(46, 14).Value = -3560.3333

$ws.Cells.Item(132, 8).Value = 3145.7144
$ws.Cells.Item(132, 9).Value = 2961.9167
$ws.Cells.Item(132, 11).Value = 8885.750100000001
$ws.Cells.Item(132, 13).Value = -6355.750100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 16900
$ws.Cells.Item(33, 10).Value = 16900
$ws.Cells.Item(33, 12).Value = 16900
$ws.Cells.Item(33, 14).Value = -17400

$ws.Cells.Item(36, 8).Value = 16900
$ws.Cells.Item(36, 10).Value = 16900
$ws.Cells.Item(36, 12).Value = 16900
$ws.Cells.Item(36, 14).Value = -17400

$ws.Cells.Item(74, 8).Value = 19587.2
$ws.Cells.Item(74, 9).Value = 18645.334
$ws.Cells.Item(74, 11).Value = 18645.334
$ws.Cells.Item(74, 13).Value = -17709.334

$ws.Cells.Item(77, 8).Value = 19587.2
$ws.Cells.Item(77, 9).Value = 18645.334
$ws.Cells.Item(77, 11).Value = 55936.00199999999
$ws.Cells.Item(77, 13).Value = -51256.00199999999

$ws.Cells.Item(100, 8).Value = 2040.5
$ws.Cells.Item(100, 9).Value = 554.2222
$ws.Cells.Item(100, 10).Value = 6499.3335
$ws.Cells.Item(100, 11).Value = 1108.4444
$ws.Cells.Item(100, 12).Value = 12998.667
$ws.Cells.Item(100, 13).Value = -567.4444000000001
$ws.Cells.Item(100, 14).Value = -14080.667
